$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "96.504.23"
Set-TextValue "E2" "  +0.80%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.590.22"
Set-TextValue "E3" "  -0.31%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.03%  "

# Row 5 - Solana
Set-TextValue "D5" "241.80"
Set-TextValue "E5" "  +1.05%  "

# Row 6 - BNB
Set-TextValue "D6" "655.56"
Set-TextValue "E6" "  +0.10%  "

# Row 7 - XRP
Set-TextValue "D7" "1.56"
Set-TextValue "E7" "  +5.87%  "

# Row 8 - Dogecoin
Set-TextValue "E8" "  -0.79%  "

# Row 9 - USDC
Set-TextValue "E9" "  +0.04%  "

# Row 10 - Cardano
Set-TextValue "E10" "  +2.86%  "

# Row 11 - LidoStakedEther
Set-TextValue "D11" "3.588.41"
Set-TextValue "E11" "  -0.29%  "

# Row 12 - Avalanche
Set-TextValue "D12" "43.25"
Set-TextValue "E12" "  -0.20%  "

# Row 13 - TRON
Set-TextValue "E13" "  +0.58%  "

# Row 14 - Toncoin
Set-TextValue "D14" "6.41"
Set-TextValue "E14" "  +0.95%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.256.30"
Set-TextValue "E15" "  -0.73%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "96.313.21"
Set-TextValue "E16" "  +0.80%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0000259"
Set-TextValue "E17" "  +0.83%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.581.24"
Set-TextValue "E18" "  -0.51%  "

# Row 19 - Polkadot
Set-TextValue "D19" "7.79"
Set-TextValue "E19" "  -2.06%  "

# Row 20 - Uniswap
Set-TextValue "D20" "12.58"
Set-TextValue "E20" "  -0.03%  "

# Row 21 - Chainlink
Set-TextValue "D21" "17.81"
Set-TextValue "E21" "  -1.62%  "

# Row 22 - Stellar
Set-TextValue "D22" "0.492"
Set-TextValue "E22" "  +0.93%  "

# Row 23 - SuiNetwork
Set-TextValue "D23" "3.46"
Set-TextValue "E23" "  -2.25%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "512.22"
Set-TextValue "E24" "  -0.10%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000202"
Set-TextValue "E25" "  +2.73%  "

# Row 26 - NEARProtocol
Set-TextValue "E26" "  +2.90%  "

# Row 27 - Litecoin
Set-TextValue "D27" "96.54"
Set-TextValue "E27" "  -0.47%  "

# Row 28 - Aptos
Set-TextValue "D28" "12.83"
Set-TextValue "E28" "  -0.02%  "

# Row 29 - WrappedeETH
Set-TextValue "D29" "3.782.57"
Set-TextValue "E29" "  -0.47%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  -7.32%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.148"
Set-TextValue "E31" "  +5.83%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "11.48"

# Row 33 - Dai
Set-TextValue "E33" "  +0.23%  "

# Row 34 - Cronos
Set-TextValue "E34" "  +3.59%  "

# Row 35 - Binance-PegBSC-USD
Set-TextValue "D35" "0.996"
Set-TextValue "E35" "  +0.01%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "31.70"
Set-TextValue "E36" "  -0.76%  "

# Row 37 - Bittensor
Set-TextValue "D37" "615.97"
Set-TextValue "E37" "  +8.52%  "

# Row 38 - swapped to RenderToken
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D38" "8.70"
Set-TextValue "E38" "  +5.12%  "

# Row 39 - swapped to PolygonEcosystemToken
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D39" "0.566"
Set-TextValue "E39" "  +0.69%  "

# Row 40 - Fetch.AI
Set-TextValue "D40" "1.62"
Set-TextValue "E40" "  +8.42%  "

# Row 41 - USDe
Set-TextValue "E41" "  +0.05%  "

# Row 42 - Kaspa
Set-TextValue "E42" "  -0.14%  "

# Row 43 - ARBITRUM
Set-TextValue "E43" "  -2.14%  "

# Row 44 - ImmutableX
Set-TextValue "E44" "  +5.81%  "

# Row 45 - Filecoin
Set-TextValue "D45" "5.71"
Set-TextValue "E45" "  -0.63%  "

# Row 46 - Stacks
Set-TextValue "E46" "  +0.94%  "

# Row 47 - EnergySwap
Set-TextValue "D47" "34.25"
Set-TextValue "E47" "  +0.47%  "

# Row 48 - WhiteBITCoin
Set-TextValue "D48" "23.53"
Set-TextValue "E48" "  -1.03%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0417"
Set-TextValue "E49" "  -0.61%  "

# Row 50 - MantraDAO
Set-TextValue "D50" "3.57"
Set-TextValue "E50" "  +3.37%  "

# Row 51 - swapped OKB -> dogwifhat
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "3.21"
Set-TextValue "E51" "  +2.40%  "
